# Regenerate merged AHB files
# - Rename the "_old"/"_new" header-label suffixes to "_FV2404"/"_FV2410"
# - Wrap the data range in an Excel Table ("Table1") with an AutoFilter
# - Freeze the header row (row 1) and select the top-left cell of the
#   scrollable area

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row labels (columns A1:J1 -> _FV2404, L1:U1 -> _FV2410,
#    K1 "diff" stays untouched).
$baseNames  = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$leftCols   = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rightCols  = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($leftCols[$i]  + "1").Value = ($baseNames[$i] + "_FV2404")
    $ws.Range($rightCols[$i] + "1").Value = ($baseNames[$i] + "_FV2410")
}

# 2) Turn the used range into a native Excel Table with an AutoFilter.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U70"), $null, 1)
$tbl.Name = "Table1"

# 3) Freeze panes above row 2 (i.e. freeze the header row) and select the
#    first scrollable cell.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
